$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New StatQuery text (Cypher query) used for Cases/Samples/Files tab rows.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Bullmastiff']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Replace the old StatQuery text in column C for the Cases/Samples/Files rows.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the view: zoom to 85% and move the selection to B4 (last edited row).
$excel.ActiveWindow.Zoom = 85
$ws.Range("B4").Select() | Out-Null
